$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "one two three"
$ws.Range("E1").Value = "date"
$ws.Range("F8").Select()
